# Updated symbol list (price/volume refresh) on the crypto sheet.
# Values are written with a leading apostrophe to force Excel to keep
# them as text (matching the original inlineStr-typed cells) instead of
# auto-converting numeric-looking strings / percentages into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'328.69"
$ws.Range("D3").Value = "'44.69"
$ws.Range("E3").Value = "'8.46%"
$ws.Range("E4").Value = "'-2.36%"
$ws.Range("D5").Value = "'0.08089"
$ws.Range("E5").Value = "'-3.63%"
$ws.Range("D6").Value = "'8.672"
$ws.Range("E6").Value = "'-1.52%"
$ws.Range("D7").Value = "'1.907"
$ws.Range("E7").Value = "'-3.79%"
$ws.Range("D8").Value = "'4.299"
$ws.Range("E8").Value = "'-4.91%"
$ws.Range("E9").Value = "'-5.92%"
$ws.Range("D10").Value = "'0.9467"
$ws.Range("E10").Value = "'2.42%"
$ws.Range("D11").Value = "'0.1185"
$ws.Range("E11").Value = "'-4.87%"
$ws.Range("D12").Value = "'0.1895"
$ws.Range("E12").Value = "'-3.19%"
$ws.Range("D13").Value = "'0.09709"
$ws.Range("E13").Value = "'3.79%"
$ws.Range("D14").Value = "'0.04060"
$ws.Range("E14").Value = "'1.42%"
$ws.Range("E15").Value = "'0.37%"
$ws.Range("D16").Value = "'0.001281"
$ws.Range("E16").Value = "'-1.77%"
$ws.Range("D17").Value = "'0.005952"
$ws.Range("E17").Value = "'-2.66%"
$ws.Range("D18").Value = "'3.576"
$ws.Range("E18").Value = "'4.17%"
$ws.Range("D19").Value = "'0.3485"
$ws.Range("E19").Value = "'-0.71%"
$ws.Range("D20").Value = "'8.520"
$ws.Range("E20").Value = "'-6.53%"
$ws.Range("D21").Value = "'0.1370"
$ws.Range("E21").Value = "'-0.41%"
$ws.Range("E22").Value = "'-1.85%"
$ws.Range("D23").Value = "'0.04355"
$ws.Range("E23").Value = "'-1.51%"
$ws.Range("D24").Value = "'0.001239"
$ws.Range("E24").Value = "'-0.49%"
$ws.Range("D25").Value = "'0.004624"
$ws.Range("E25").Value = "'6.39%"
$ws.Range("E26").Value = "'2.94%"
$ws.Range("D27").Value = "'0.0003993"
$ws.Range("E27").Value = "'-0.35%"
$ws.Range("E39").Value = "'-4.74%"
$ws.Range("D40").Value = "'0.05595"
$ws.Range("E40").Value = "'1.33%"
$ws.Range("E41").Value = "'25.91%"
$ws.Range("D42").Value = "'0.007681"
$ws.Range("E42").Value = "'-3.26%"
$ws.Range("D43").Value = "'0.1397"
$ws.Range("E43").Value = "'-2.27%"
$ws.Range("E44").Value = "'-3.12%"
$ws.Range("D45").Value = "'0.009445"
$ws.Range("E45").Value = "'-6.62%"
$ws.Range("D46").Value = "'0.00007082"
$ws.Range("E46").Value = "'-1.37%"
$ws.Range("E47").Value = "'-0.35%"
$ws.Range("D48").Value = "'0.003446"
$ws.Range("E48").Value = "'-0.57%"
$ws.Range("D49").Value = "'0.002272"
$ws.Range("E49").Value = "'-0.69%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'-0.35%"
$ws.Range("E51").Value = "'-0.35%"
